$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A14,L15:Q15")
Write-Host $r.Address
Write-Host $r.Areas.Count
